$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 1872.39250303085
$ws.Range("E4").Value = 1878.0448866
$ws.Range("F4").Value = 1883.62592328678
$ws.Range("G4").Value = 1889.08350979391
$ws.Range("H4").Value = 1894.77691976033
$ws.Range("I4").Value = 1904.67975232033
$ws.Range("J4").Value = 1914.95303263777
$ws.Range("K4").Value = 1925.96441890162
# Row 5
$ws.Range("D5").Value = 3175.98175838004
$ws.Range("E5").Value = 3206.09574740905
$ws.Range("F5").Value = 3232.36567083857
$ws.Range("G5").Value = 3261.90611278909
$ws.Range("H5").Value = 3292.54940462156
$ws.Range("I5").Value = 3324.0646831367
$ws.Range("J5").Value = 3353.44808520622
$ws.Range("K5").Value = 3383.38513261417
# Row 10
$ws.Range("D10").Value = 674.033916028452
$ws.Range("E10").Value = 684.10547609577
$ws.Range("F10").Value = 700.568454769733
$ws.Range("G10").Value = 718.736868736257
$ws.Range("H10").Value = 733.964067457659
$ws.Range("I10").Value = 745.270944386688
$ws.Range("J10").Value = 756.730817045884
$ws.Range("K10").Value = 768.345755652654
# Row 11
$ws.Range("D11").Value = 980.088833465518
$ws.Range("E11").Value = 994.733532082047
$ws.Range("F11").Value = 1018.67176601989
$ws.Range("G11").Value = 1045.08981298597
$ws.Range("H11").Value = 1067.23114308364
$ws.Range("I11").Value = 1083.67207217635
$ws.Range("J11").Value = 1100.33546693903
$ws.Range("K11").Value = 1117.22433760144
# Row 12
$ws.Range("D12").Value = 1116.08378227538
$ws.Range("E12").Value = 1130.29719654266
$ws.Range("F12").Value = 1144.69161975241
$ws.Range("G12").Value = 1159.26935706766
$ws.Range("H12").Value = 1176.61781913306
$ws.Range("I12").Value = 1194.22590087546
$ws.Range("J12").Value = 1212.09748750246
$ws.Range("K12").Value = 1230.23652236378
# Row 19
$ws.Range("D19").Value = 4506.17006477277
$ws.Range("E19").Value = 4663.80185798619
$ws.Range("F19").Value = 4687.3360681327
$ws.Range("G19").Value = 4711.08401475736
$ws.Range("H19").Value = 4735.045843276
$ws.Range("I19").Value = 5040.33717271295
$ws.Range("J19").Value = 5066.12244483538
$ws.Range("K19").Value = 5092.12908886204
# Row 20
$ws.Range("D20").Value = 2498.70051302068
$ws.Range("E20").Value = 2529.15816810511
$ws.Range("F20").Value = 2560.6019168812
$ws.Range("G20").Value = 2587.60624358038
$ws.Range("H20").Value = 2614.20733389807
$ws.Range("I20").Value = 2638.67223947288
$ws.Range("J20").Value = 2664.12111126134
$ws.Range("K20").Value = 2687.35532427355
# Row 21
$ws.Range("D21").Value = 509.332755158353
$ws.Range("E21").Value = 526.099889428512
$ws.Range("F21").Value = 543.418994465887
$ws.Range("G21").Value = 561.308241039732
$ws.Range("H21").Value = 557.839987860135
$ws.Range("I21").Value = 554.393164581683
$ws.Range("J21").Value = 550.967638791707
$ws.Range("K21").Value = 547.563278895698
# Row 22
$ws.Range("D22").Value = 166.108108108108
$ws.Range("E22").Value = 166.550368550369
$ws.Range("F22").Value = 166.481572481572
$ws.Range("G22").Value = 166.624078624079
$ws.Range("H22").Value = 166.540540540541
$ws.Range("I22").Value = 166.742014742015
$ws.Range("J22").Value = 167.071253071253
$ws.Range("K22").Value = 166.904176904177
# Row 24
$ws.Range("C24").Value = 79.6290028317487
$ws.Range("D24").Value = 71.1791672183253
$ws.Range("E24").Value = 73.9806362924226
$ws.Range("F24").Value = 88.3944428040623
$ws.Range("G24").Value = 103.083510693309
$ws.Range("H24").Value = 109.658761987664
$ws.Range("I24").Value = 134.709036355557
$ws.Range("J24").Value = 155.206844624839
$ws.Range("K24").Value = 162.820344255939
